$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163, shifting existing rows 163:204 down to 164:205
$ws.Rows(163).Insert()

# Populate the new row 163 with its data
$ws.Cells.Item(163, 1).Value = 5
$ws.Cells.Item(163, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(163, 3).Value = 'Maule'
$ws.Cells.Item(163, 4).Value = '2022-03-22'
$ws.Cells.Item(163, 5).Value = 7
$ws.Cells.Item(163, 6).Value = 100112008
$ws.Cells.Item(163, 7).Value = 'Coliflor'
$ws.Cells.Item(163, 8).Value = 'Sin especificar'
$ws.Cells.Item(163, 9).Value = 'Primera'
$ws.Cells.Item(163, 10).Value = 2000
$ws.Cells.Item(163, 11).Value = 1000
$ws.Cells.Item(163, 12).Value = 1000
$ws.Cells.Item(163, 13).Value = 1000
$ws.Cells.Item(163, 14).Value = '$/unidad'
$ws.Cells.Item(163, 15).Value = 'Región del Maule'
$ws.Cells.Item(163, 16).Value = 1000
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = 'Hortaliza'
